$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in this sheet hold text (prices use "." as a thousands
# separator, percentages keep padding spaces), so force the Text number
# format before writing to stop Excel from auto-coercing the string into a
# numeric value (which would also lose the original formatting/precision).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.671.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.322.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.30%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('B5').NumberFormat = '@'
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').NumberFormat = '@'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '270.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.27'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.99%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.04'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.08'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.71%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.670.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.69'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.854'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +8.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.328.22'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.604.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.35'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +6.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.94'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.37'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.40'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +9.69%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.53'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.39'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.47'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.29'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.24'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.50'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '172.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0898'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.48'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.98%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.70%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.15%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.85%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.34'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.33'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.63%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +10.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.34'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +17.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.09'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +8.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.85'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.36'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.21'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.551.44'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.182'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +13.63%  '
